# Scheduled-runner update: refresh currentAveragePrice* / Leve profit figures
# across the ALC / BSM / CRP / CUL / LTW / WVR sheets (Ultima_Profits tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4800
$ws.Range("I69").Value = 6750
$ws.Range("K69").Value = 20250
$ws.Range("M69").Value = -19376

$ws.Range("H72").Value = 4800
$ws.Range("I72").Value = 6750
$ws.Range("K72").Value = 60750
$ws.Range("M72").Value = -56382

$ws.Range("H80").Value = 2452.125
$ws.Range("I80").Value = 2336.1667
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 7008.500100000001
$ws.Range("L80").Value = 8400
$ws.Range("M80").Value = -6010.500100000001
$ws.Range("N80").Value = -10396

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H82").Value = 3379.6667
$ws.Range("I82").Value = 546
$ws.Range("J82").Value = 9047
$ws.Range("K82").Value = 1638
$ws.Range("L82").Value = 27141
$ws.Range("M82").Value = -1232
$ws.Range("N82").Value = -27953

$ws.Range("H83").Value = 2452.125
$ws.Range("I83").Value = 2336.1667
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 21025.5003
$ws.Range("L83").Value = 25200
$ws.Range("M83").Value = -16033.5003
$ws.Range("N83").Value = -35184

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H85").Value = 3379.6667
$ws.Range("I85").Value = 546
$ws.Range("J85").Value = 9047
$ws.Range("K85").Value = 1638
$ws.Range("L85").Value = 27141
$ws.Range("M85").Value = -234
$ws.Range("N85").Value = -29949

$ws.Range("H86").Value = 2926.8518
$ws.Range("I86").Value = 2333
$ws.Range("J86").Value = 5539.8
$ws.Range("K86").Value = 2333
$ws.Range("L86").Value = 5539.8
$ws.Range("M86").Value = -1210
$ws.Range("N86").Value = -7785.8

$ws.Range("H87").Value = 57020.668
$ws.Range("J87").Value = 57020.668
$ws.Range("L87").Value = 57020.668
$ws.Range("N87").Value = -59516.668

$ws.Range("H88").Value = 14500
$ws.Range("I88").Value = 14500
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 14500
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -14094
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 2926.8518
$ws.Range("I89").Value = 2333
$ws.Range("J89").Value = 5539.8
$ws.Range("K89").Value = 11665
$ws.Range("L89").Value = 27699
$ws.Range("M89").Value = -6049
$ws.Range("N89").Value = -38931

$ws.Range("H90").Value = 57020.668
$ws.Range("J90").Value = 57020.668
$ws.Range("L90").Value = 171062.004
$ws.Range("N90").Value = -183542.004

$ws.Range("H91").Value = 14500
$ws.Range("I91").Value = 14500
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 14500
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -13096
$ws.Range("N91").ClearContents()

$ws.Range("H94").Value = 4000
$ws.Range("I94").Value = 4000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3549
$ws.Range("N94").ClearContents()

$ws.Range("H98").Value = 2225
$ws.Range("I98").Value = 2225
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2225
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -727
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 2225
$ws.Range("I122").Value = 2225
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6675
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4225
$ws.Range("N122").Value = -4225

$ws.Range("H138").Value = 6412085.5
$ws.Range("I138").Value = 1424.1765
$ws.Range("J138").Value = 11365778
$ws.Range("K138").Value = 4272.529500000001
$ws.Range("L138").Value = 34097334
$ws.Range("M138").Value = 867.4704999999994
$ws.Range("N138").Value = -34107614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 13500
$ws.Range("J40").Value = 13500
$ws.Range("L40").Value = 13500
$ws.Range("N40").Value = -14030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6947782
$ws.Range("I31").Value = 3264.8262
$ws.Range("J31").Value = 166671680
$ws.Range("K31").Value = 3264.8262
$ws.Range("L31").Value = 166671680
$ws.Range("M31").Value = -2969.8262
$ws.Range("N31").Value = -166672270

$ws.Range("H34").Value = 6947782
$ws.Range("I34").Value = 3264.8262
$ws.Range("J34").Value = 166671680
$ws.Range("K34").Value = 3264.8262
$ws.Range("L34").Value = 166671680
$ws.Range("M34").Value = -3062.8262
$ws.Range("N34").Value = -166672084

$ws.Range("H129").Value = 41001.5
$ws.Range("J129").Value = 42456.184
$ws.Range("L129").Value = 42456.184
$ws.Range("N129").Value = -52456.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 925
$ws.Range("I35").Value = 50
$ws.Range("J35").Value = 1800
$ws.Range("K35").Value = 150
$ws.Range("L35").Value = 5400
$ws.Range("M35").Value = 138
$ws.Range("N35").Value = -5976

$ws.Range("H39").Value = 520.5789
$ws.Range("I39").Value = 400
$ws.Range("J39").Value = 523.8378
$ws.Range("K39").Value = 1200
$ws.Range("L39").Value = 1571.5134
$ws.Range("M39").Value = -906
$ws.Range("N39").Value = -2159.5134

$ws.Range("H118").Value = 3069.3076
$ws.Range("I118").Value = 2739.25
$ws.Range("K118").Value = 8217.75
$ws.Range("M118").Value = -6974.75

$ws.Range("H129").Value = 2277.1785
$ws.Range("I129").Value = 890.9091
$ws.Range("J129").Value = 3174.1765
$ws.Range("K129").Value = 2672.7273
$ws.Range("L129").Value = 9522.529500000001
$ws.Range("M129").Value = 2327.2727
$ws.Range("N129").Value = -19522.5295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6708.231
$ws.Range("I40").Value = 14250
$ws.Range("J40").Value = 3356.3333
$ws.Range("K40").Value = 14250
$ws.Range("L40").Value = 3356.3333
$ws.Range("M40").Value = -14114
$ws.Range("N40").Value = -3628.3333

$ws.Range("H132").Value = 10424743
$ws.Range("I132").Value = 5145.933
$ws.Range("J132").Value = 27790740
$ws.Range("K132").Value = 15437.799
$ws.Range("L132").Value = 83372220
$ws.Range("M132").Value = -12907.799
$ws.Range("N132").Value = -83377280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1457.4138
$ws.Range("I132").Value = 1032.1765
$ws.Range("J132").Value = 2059.8333
$ws.Range("K132").Value = 3096.5295
$ws.Range("L132").Value = 6179.499899999999
$ws.Range("M132").Value = -566.5295000000001
$ws.Range("N132").Value = -11239.4999
